$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-20 20:55:17"
$wsZh.Range("H2").Value = "2016-03-20 20:55:37"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-20 20:55:20"
$wsDe.Range("H2").Value = "2016-03-20 20:55:43"
